# Edu sheet ("Get consistent values over time"):
# Expand the two small "C01 / C02 percent" legend blocks (rows 24-29) into
# full six-row blocks that mirror the Pct HS / Pct BA / Pop25+LessHS:pov /
# HS:pov / College:pov / BA:pov rows, each tagged with its C01/C02 column
# label, and move the "2015-2017" block down to make room (leaving row 31
# blank), matching the analogous PovFam sheet layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Edu")

# --- "2010-2014" block (C01) -------------------------------------------------
$ws.Range("A24").Value = "2010-2014"

$ws.Range("A25").Value = "Pct HS"
$ws.Range("B25").Value = 14
$ws.Range("C25").Value = "percent"
$ws.Range("D25").Value = "C01"

$ws.Range("A26").Value = "Pct BA"
$ws.Range("B26").Value = 15
$ws.Range("C26").Value = "percent"
$ws.Range("D26").Value = "C01"

$ws.Range("A27").Value = "Pop 25+ Less HS: pov"
$ws.Range("B27").Value = 28
$ws.Range("C27").Value = "percent"
$ws.Range("D27").Value = "C01"

$ws.Range("A28").Value = "HS: pov"
$ws.Range("B28").Value = 29
$ws.Range("C28").Value = "percent"
$ws.Range("D28").Value = "C01"

$ws.Range("A29").Value = "College: pov"
$ws.Range("B29").Value = 30
$ws.Range("C29").Value = "percent"
$ws.Range("D29").Value = "C01"

$ws.Range("A30").Value = "BA: pov"
$ws.Range("B30").Value = 31
$ws.Range("C30").Value = "percent"
$ws.Range("D30").Value = "C01"

# Row 31 intentionally left blank (gap between the two blocks).

# --- "2015-2017" block (C02), now starting at row 32 ------------------------
$ws.Range("A32").Value = "2015-2017"

$ws.Range("A33").Value = "Pct HS"
$ws.Range("B33").Value = 14
$ws.Range("C33").Value = "percent"
$ws.Range("D33").Value = "C02"

$ws.Range("A34").Value = "Pct BA"
$ws.Range("B34").Value = 15
$ws.Range("C34").Value = "percent"
$ws.Range("D34").Value = "C02"

$ws.Range("A35").Value = "Pop 25+ Less HS: pov"
$ws.Range("B35").Value = 55
$ws.Range("C35").Value = "percent"
$ws.Range("D35").Value = "C02"

$ws.Range("A36").Value = "HS: pov"
$ws.Range("B36").Value = 56
$ws.Range("C36").Value = "percent"
$ws.Range("D36").Value = "C02"

$ws.Range("A37").Value = "College: pov"
$ws.Range("B37").Value = 57
$ws.Range("C37").Value = "percent"
$ws.Range("D37").Value = "C02"

$ws.Range("A38").Value = "BA: pov"
$ws.Range("B38").Value = 58
$ws.Range("C38").Value = "percent"
$ws.Range("D38").Value = "C02"

# Match the author's final selection/scroll state on the Edu tab.
$ws.Activate()
$ws.Range("F34").Select()
